$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from the last existing data row (255) down to the new rows (256:269)
$ws.Range("A255:D255").Copy($ws.Range("A256:D269"))

$data = @(
    @(256, 44330, 1, 9, 144.2076590290018),
    @(257, 44331, 0, 9, 144.2076590290018),
    @(258, 44332, 0, 9, 144.2076590290018),
    @(259, 44333, 0, 4, 64.09229290177856),
    @(260, 44334, 1, 4, 64.09229290177856),
    @(261, 44335, 0, 3, 48.06921967633392),
    @(262, 44336, 3, 5, 80.11536612722321),
    @(263, 44337, 0, 4, 64.09229290177856),
    @(264, 44338, 0, 4, 64.09229290177856),
    @(265, 44339, 0, 4, 64.09229290177856),
    @(266, 44340, 0, 4, 64.09229290177856),
    @(267, 44341, 0, 3, 48.06921967633392),
    @(268, 44342, 0, 3, 48.06921967633392),
    @(269, 44343, 0, 0, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}

Write-Host "Done updating rows 256-269"
